$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New cell A13 holding the shared-string value "Actif"
$ws.Range("A13").Value = "Actif"

# New row 19 with a formula that mirrors the A15/A16 pattern but tests A1
# and returns custom "true"/"false" sentences.
$ws.Range("A19").Formula = '=IF(A1=TRUE,"Votre phrase personnalisée si vrai","Votre phrase personnalisée si faux")'

# Move the selection to the newly added cell, matching the saved view state.
$ws.Range("A19").Select()
